# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 11:47"

# Finlandia (row 52): refresh case counts
$ws.Cells.Item(52, 2).Value = 4576
$ws.Cells.Item(52, 3).Value = 101
$ws.Cells.Item(52, 5).Value = 1890

# Costa de Marfil (row 86): refresh case counts
$ws.Cells.Item(86, 2).Value = 1111
$ws.Cells.Item(86, 3).Value = 34
$ws.Cells.Item(86, 4).Value = 449
$ws.Cells.Item(86, 5).Value = 648

# Madagascar / Etiopia (rows 140-141): swap alphabetical order and refresh values
$ws.Cells.Item(140, 1).Value = "Etiopia"
$ws.Cells.Item(140, 2).Value = 123
$ws.Cells.Item(140, 3).Value = 1
$ws.Cells.Item(140, 4).Value = 41
$ws.Cells.Item(140, 5).Value = 79
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 3

$ws.Cells.Item(141, 1).Value = "Madagascar"
$ws.Cells.Item(141, 2).Value = 123
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 62
$ws.Cells.Item(141, 5).Value = 61
$ws.Cells.Item(141, 6).Value = 1
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 0
